# UND-EFF figure tweak: nudge the "Estimated gap" annotation
# (arrow connector, horizontal line connector and its textbox label)
# 30890 EMU to the right, matching the author's re-positioning edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU delta to apply along the x-axis (matches the OOXML diff exactly).
$deltaEmu = 30890
$emuPerPt = 12700.0

# Shape.Left/Top round-trip through a single-precision (float) COM
# property, so nudge the computed point value by half an EMU before
# assigning it; this compensates for the float32 truncation that would
# otherwise occasionally leave the stored EMU one unit short.
$halfEmuPt = 0.5 / $emuPerPt

# Target shape ids (drawing.xml creationIds in parentheses for reference):
#   149 - Straight Arrow Connector 148 ({0D7C4CEF-49CD-4693-B536-8214DE9B278B})
#   150 - Straight Connector 149        ({E6968A13-247B-406E-A618-8B4828738B85})
#   151 - TextBox 150 "Estimated gap"   ({15FADE80-A39C-4E9F-85C7-24A19456D068})
$targetIds = @(149, 150, 151)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    $id = [int]$sh.Id
    if ($targetIds -contains $id) {
        $currentEmu = [math]::Round($sh.Left * $emuPerPt)
        $newEmu = $currentEmu + $deltaEmu
        $sh.Left = ($newEmu / $emuPerPt) + $halfEmuPt
    }
}
